$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete entire row 159 (the "「試験」" entry) - all rows below shift up by one.
$ws.Rows.Item(159).Delete()
